# Fix service query, add in refrigerator status form, fix formatting dates in list view
#
# The "choices" sheet has a `service_priority_list` choice list (rows 59-62)
# used by the refrigerator status form. The "not_applicable" option had been
# appended at the end of the list; move it to the front (right after the
# list name) so it appears first in the form, ahead of low/medium/high.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

$firstRow = 59
$lastRow  = 62
$count    = $lastRow - $firstRow + 1

# Capture the data_value / label / label_es (columns B:D) for the rows that
# make up the service_priority_list choice list.
$rows = New-Object 'object[]' $count
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $idx = $r - $firstRow
    $rows[$idx] = @($ws.Cells.Item($r, 2).Value2, $ws.Cells.Item($r, 3).Value2, $ws.Cells.Item($r, 4).Value2)
}

# Reorder: move the last entry (not_applicable) to the front of the list;
# low/medium/high keep their relative order but shift down by one row.
$reordered = New-Object 'object[]' $count
$reordered[0] = $rows[$count - 1]
for ($i = 1; $i -lt $count; $i++) {
    $reordered[$i] = $rows[$i - 1]
}

$r = $firstRow
for ($i = 0; $i -lt $count; $i++) {
    $row = $reordered[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# Reflect the resulting cursor position in the sheet view.
$ws.Activate()
[void]$ws.Range("B71").Select()
